# Apply weekly report data refresh:
#  - Update the "Report Generated On" timestamp
#  - Populate Total Billed Amount / line item pricing that were previously 0
#  - Clear the now-unused Scope ID # value

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Report generation timestamp (D5)
$ws.Range("D5").Value = "Report Generated On: 08/26/2025 10:00 AM"

# Total Billed Amount (C8)
$ws.Range("C8").Value = 2000.67

# Scope ID # (G10) is no longer populated
$ws.Range("G10").Value = ""

# Per-line-item pricing (column H) for the detail rows
$ws.Range("H16").Value = 94.17
$ws.Range("H17").Value = 94.17
$ws.Range("H18").Value = 478.55
$ws.Range("H19").Value = 94.17
$ws.Range("H20").Value = 94.17
$ws.Range("H21").Value = 478.55
$ws.Range("H22").Value = 94.17
$ws.Range("H23").Value = 94.17
$ws.Range("H24").Value = 478.55

# TOTAL pricing (H25)
$ws.Range("H25").Value = 2000.67
